# Update "Paises" (countries) COVID data sheet and provincias-style refresh.
# - Refresh the "last updated" timestamp (16:05 -> 16:35)
# - Update Azerbaiyan / Cuba / Yemen rows with refreshed figures, which pushes them
#   up past the countries they previously trailed (Camerun/Grecia, Sudan, Barbados/
#   Liechtenstein/San Martin (Parte Holandesa) respectively), since the sheet is kept
#   sorted by total cases (column B) descending.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1) Refresh timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 14 de Mayo de 2020 a las 16:35"

# 2) Row 4 (Estados Unidos) refreshed totals
$ws.Range("B4").Value = 1433613
$ws.Range("C4").Value = 3265
$ws.Range("E4").Value = 1037856
$ws.Range("G4").Value = 145
$ws.Range("H4").Value = 85342

# 3) Row 20 (Paises Bajos) refreshed totals
$ws.Range("F20").Value = 410

# 4) Row 54 (Argentina) refreshed totals
$ws.Range("D54").Value = 2385
$ws.Range("E54").Value = 4150
$ws.Range("G54").Value = 15
$ws.Range("H54").Value = 344

# 5) Azerbaiyan moves ahead of Camerun & Grecia (rows 72-74) with fresh data,
#    remaining two rows simply shift down keeping their previous numbers.
$ws.Range("A72").Value = "Azerbaiyan"
$ws.Range("B72").Value = 2879
$ws.Range("C72").Value = 121
$ws.Range("D72").Value = 1833
$ws.Range("E72").Value = 1011
$ws.Range("F72").Value = 29
$ws.Range("G72").Value = 0
$ws.Range("H72").Value = 35

$ws.Range("A73").Value = "Camerun"
$ws.Range("B73").Value = 2800
$ws.Range("C73").Value = 0
$ws.Range("D73").Value = 1543
$ws.Range("E73").Value = 1121
$ws.Range("F73").Value = 28
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 136

$ws.Range("A74").Value = "Grecia"
$ws.Range("B74").Value = 2760
$ws.Range("C74").Value = 0
$ws.Range("D74").Value = 1374
$ws.Range("E74").Value = 1231
$ws.Range("F74").Value = 28
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 155

# 6) Cuba moves ahead of Sudan (rows 83-84) with fresh data.
$ws.Range("A83").Value = "Cuba"
$ws.Range("B83").Value = 1830
$ws.Range("C83").Value = 20
$ws.Range("D83").Value = 1383
$ws.Range("E83").Value = 368
$ws.Range("F83").Value = 9
$ws.Range("G83").Value = 0
$ws.Range("H83").Value = 79

$ws.Range("A84").Value = "Sudan"
$ws.Range("B84").Value = 1818
$ws.Range("C84").Value = 0
$ws.Range("D84").Value = 198
$ws.Range("E84").Value = 1530
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 90

# 7) Yemen moves ahead of Barbados, Liechtenstein & San Martin (Parte Holandesa)
#    (rows 168-171) with fresh data, others shift down keeping their previous numbers.
$ws.Range("A168").Value = "Yemen"
$ws.Range("B168").Value = 85
$ws.Range("C168").Value = 15
$ws.Range("D168").Value = 1
$ws.Range("E168").Value = 72
$ws.Range("F168").Value = 0
$ws.Range("G168").Value = 0
$ws.Range("H168").Value = 12

$ws.Range("A169").Value = "Barbados"
$ws.Range("B169").Value = 85
$ws.Range("C169").Value = 0
$ws.Range("D169").Value = 65
$ws.Range("E169").Value = 13
$ws.Range("F169").Value = 4
$ws.Range("G169").Value = 0
$ws.Range("H169").Value = 7

$ws.Range("A170").Value = "Liechtenstein"
$ws.Range("B170").Value = 82
$ws.Range("C170").Value = 0
$ws.Range("D170").Value = 55
$ws.Range("E170").Value = 26
$ws.Range("F170").Value = 0
$ws.Range("G170").Value = 0
$ws.Range("H170").Value = 1

$ws.Range("A171").Value = "San Martin (Parte Holandesa)"
$ws.Range("B171").Value = 76
$ws.Range("C171").Value = 0
$ws.Range("D171").Value = 46
$ws.Range("E171").Value = 15
$ws.Range("F171").Value = 7
$ws.Range("G171").Value = 0
$ws.Range("H171").Value = 15
